$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.047605037959895
$ws.Range("D2").Value = 1.054620021756721
$ws.Range("E2").Value = 1.055423923710415
$ws.Range("F2").Value = 1.067056793130225
$ws.Range("I2").Value = 1.049329951504624
$ws.Range("J2").Value = 1.052653217580459
$ws.Range("K2").Value = 1.05736273107179
$ws.Range("L2").Value = 1.058164420647572
$ws.Range("M2").Value = 1.069765679598073
$ws.Range("N2").Value = 1.054148105890435
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.048453455319179
$ws.Range("D3").Value = 1.055281909053358
$ws.Range("E3").Value = 1.056165507218695
$ws.Range("F3").Value = 1.067848162409601
$ws.Range("I3").Value = 1.049566141559956
$ws.Range("J3").Value = 1.053150672347866
$ws.Range("K3").Value = 1.057838284587375
$ws.Range("L3").Value = 1.058719625999357
$ws.Range("M3").Value = 1.070372821304469
$ws.Range("N3").Value = 1.05464626710067
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.049003089843431
$ws.Range("D4").Value = 1.055710737526487
$ws.Range("E4").Value = 1.056646291962639
$ws.Range("F4").Value = 1.068361177332577
$ws.Range("I4").Value = 1.049718015016554
$ws.Range("J4").Value = 1.053472528810654
$ws.Range("K4").Value = 1.058145856886312
$ws.Range("L4").Value = 1.059079141626038
$ws.Range("M4").Value = 1.070765967523251
$ws.Range("N4").Value = 1.054968580636554
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.049234310538335
$ws.Range("D5").Value = 1.055891145229193
$ws.Range("E5").Value = 1.056848634776689
$ws.Range("F5").Value = 1.068577073310085
$ws.Range("I5").Value = 1.049781632586371
$ws.Range("J5").Value = 1.053607828753695
$ws.Range("K5").Value = 1.058275124783592
$ws.Range("L5").Value = 1.059230342918948
$ws.Range("M5").Value = 1.070931312783512
$ws.Range("N5").Value = 1.055104072721033
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.049273142511926
$ws.Range("D6").Value = 1.055921443941963
$ws.Range("E6").Value = 1.056882621924788
$ws.Range("F6").Value = 1.06861333630183
$ws.Range("I6").Value = 1.049792300742414
$ws.Range("J6").Value = 1.05363054567175
$ws.Range("K6").Value = 1.058296827311455
$ws.Range("L6").Value = 1.059255733811063
$ws.Range("M6").Value = 1.070959078837886
$ws.Range("N6").Value = 1.055126821899717
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.049006178821335
$ws.Range("D7").Value = 1.055713147641536
$ws.Range("E7").Value = 1.05664899481158
$ws.Range("F7").Value = 1.06836406126555
$ws.Range("I7").Value = 1.049718865982564
$ws.Range("J7").Value = 1.053474336730229
$ws.Range("K7").Value = 1.05814758431146
$ws.Range("L7").Value = 1.059081161747579
$ws.Range("M7").Value = 1.0707681766151
$ws.Range("N7").Value = 1.054970391123582
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.047891628609269
$ws.Range("D8").Value = 1.054843596117842
$ws.Range("E8").Value = 1.055674351503391
$ws.Range("F8").Value = 1.067324042994726
$ws.Range("I8").Value = 1.049409970743357
$ws.Range("J8").Value = 1.052821340171554
$ws.Range("K8").Value = 1.057523475327247
$ws.Range("L8").Value = 1.058352000042691
$ws.Range("M8").Value = 1.069970806005905
$ws.Range("N8").Value = 1.054316467234893
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.045932715347003
$ws.Range("D9").Value = 1.053315575179142
$ws.Range("E9").Value = 1.053964116325521
$ws.Range("F9").Value = 1.06549873176675
$ws.Range("I9").Value = 1.048858363755912
$ws.Range("J9").Value = 1.05167050529167
$ws.Range("K9").Value = 1.056422683024865
$ws.Range("L9").Value = 1.057069185307022
$ws.Range("M9").Value = 1.068567991215311
$ws.Range("N9").Value = 1.053163998037478
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.044630280501516
$ws.Range("D10").Value = 1.052299853516218
$ws.Range("E10").Value = 1.052828920286111
$ws.Range("F10").Value = 1.064286907266578
$ws.Range("I10").Value = 1.048485772422544
$ws.Range("J10").Value = 1.050903245632558
$ws.Range("K10").Value = 1.055688207941997
$ws.Range("L10").Value = 1.056215449350748
$ws.Range("M10").Value = 1.067634393812666
$ws.Range("N10").Value = 1.052395648781643
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.044067164935121
$ws.Range("D11").Value = 1.051860759536423
$ws.Range("E11").Value = 1.05233856829752
$ws.Range("F11").Value = 1.063763396705775
$ws.Range("I11").Value = 1.048323295811278
$ws.Range("J11").Value = 1.050571021735996
$ws.Range("K11").Value = 1.055370043619179
$ws.Range("L11").Value = 1.055846140406962
$ws.Range("M11").Value = 1.067230538757094
$ws.Range("N11").Value = 1.052062953089042
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.043858127742896
$ws.Range("D12").Value = 1.051697770538872
$ws.Range("E12").Value = 1.05215661126317
$ws.Range("F12").Value = 1.06356912652011
$ws.Range("I12").Value = 1.048262774021272
$ws.Range("J12").Value = 1.050447621016955
$ws.Range("K12").Value = 1.055251844795637
$ws.Range("L12").Value = 1.055709018861907
$ws.Range("M12").Value = 1.067080590573692
$ws.Range("N12").Value = 1.051939377126825
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.043902961065254
$ws.Range("D13").Value = 1.051732727222048
$ws.Range("E13").Value = 1.052195633419194
$ws.Range("F13").Value = 1.063610789732493
$ws.Range("I13").Value = 1.04827576386352
$ws.Range("J13").Value = 1.050474090776816
$ws.Range("K13").Value = 1.055277199656448
$ws.Range("L13").Value = 1.05573842932974
$ws.Range("M13").Value = 1.067112752155514
$ws.Range("N13").Value = 1.051965884476782
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.044049883214148
$ws.Range("D14").Value = 1.051847284559188
$ws.Range("E14").Value = 1.052323523958025
$ws.Range("F14").Value = 1.063747334487
$ws.Range("I14").Value = 1.048318296539427
$ws.Range("J14").Value = 1.05056082134495
$ws.Range("K14").Value = 1.055360273635561
$ws.Range("L14").Value = 1.055834804741323
$ws.Range("M14").Value = 1.067218142727992
$ws.Range("N14").Value = 1.05205273821227
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.044140423944094
$ws.Range("D15").Value = 1.05191788175456
$ws.Range("E15").Value = 1.05240234565183
$ws.Range("F15").Value = 1.063831488784879
$ws.Range("I15").Value = 1.048344479725137
$ws.Range("J15").Value = 1.050614259219242
$ws.Range("K15").Value = 1.05541145585495
$ws.Range("L15").Value = 1.05589419231039
$ws.Range("M15").Value = 1.06728308554555
$ws.Range("N15").Value = 1.052106251974474
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.044667670624284
$ws.Range("D16").Value = 1.052329010069765
$ws.Range("E16").Value = 1.052861488683202
$ws.Range("F16").Value = 1.064321676762569
$ws.Range("I16").Value = 1.048496531469464
$ws.Range("J16").Value = 1.050925294449474
$ws.Range("K16").Value = 1.055709320800925
$ws.Range("L16").Value = 1.056239966973679
$ws.Range("M16").Value = 1.067661204885655
$ws.Range("N16").Value = 1.052417728910409
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.044998626588356
$ws.Range("D17").Value = 1.052587094128919
$ws.Range("E17").Value = 1.053149818381235
$ws.Range("F17").Value = 1.064629486114503
$ws.Range("I17").Value = 1.048591604334232
$ws.Range("J17").Value = 1.051120400807331
$ws.Range("K17").Value = 1.055896129300853
$ws.Range("L17").Value = 1.056456960892605
$ws.Range("M17").Value = 1.067898497090367
$ws.Range("N17").Value = 1.052613112341673
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.0451917492423
$ws.Range("D18").Value = 1.052737699559064
$ws.Range("E18").Value = 1.053318111423947
$ws.Range("F18").Value = 1.064809143430883
$ws.Range("I18").Value = 1.048646948474314
$ws.Range("J18").Value = 1.05123420339896
$ws.Range("K18").Value = 1.056005078733122
$ws.Range("L18").Value = 1.056583564829122
$ws.Range("M18").Value = 1.068036943984394
$ws.Range("N18").Value = 1.052727076546035
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.045257612872285
$ws.Range("D19").Value = 1.052789063827775
$ws.Range("E19").Value = 1.05337551448485
$ws.Range("F19").Value = 1.064870421772018
$ws.Range("I19").Value = 1.048665800669254
$ws.Range("J19").Value = 1.051273007161147
$ws.Range("K19").Value = 1.056042225460405
$ws.Range("L19").Value = 1.056626739400025
$ws.Range("M19").Value = 1.068084157254639
$ws.Range("N19").Value = 1.052765935414015
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.044963109688914
$ws.Range("D20").Value = 1.052559396963737
$ws.Range("E20").Value = 1.053118871402906
$ws.Range("F20").Value = 1.06459644893351
$ws.Range("I20").Value = 1.04858141531353
$ws.Range("J20").Value = 1.051099467695331
$ws.Range("K20").Value = 1.055876087831375
$ws.Range("L20").Value = 1.056433675862008
$ws.Range("M20").Value = 1.067873033902238
$ws.Range("N20").Value = 1.052592149502252
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.044006614728469
$ws.Range("D21").Value = 1.051813547225665
$ws.Range("E21").Value = 1.052285858355371
$ws.Range("F21").Value = 1.06370712033831
$ws.Range("I21").Value = 1.048305776432991
$ws.Range("J21").Value = 1.050535281287602
$ws.Range("K21").Value = 1.055335810916985
$ws.Range("L21").Value = 1.055806423020165
$ws.Range("M21").Value = 1.067187106118392
$ws.Range("N21").Value = 1.052027161885111
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.043405975593545
$ws.Range("D22").Value = 1.051345239566188
$ws.Range("E22").Value = 1.05176316120594
$ws.Range("F22").Value = 1.063149035201203
$ws.Range("I22").Value = 1.048131483886528
$ws.Range("J22").Value = 1.050180567298524
$ws.Range("K22").Value = 1.054996011437185
$ws.Range("L22").Value = 1.055412370392735
$ws.Range("M22").Value = 1.066756193115514
$ws.Range("N22").Value = 1.051671944161483
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.043724314367346
$ws.Range("D23").Value = 1.051593437254302
$ws.Range("E23").Value = 1.052040152566343
$ws.Range("F23").Value = 1.063444784488413
$ws.Range("I23").Value = 1.048223972978236
$ws.Range("J23").Value = 1.050368606273707
$ws.Range("K23").Value = 1.05517615515565
$ws.Range("L23").Value = 1.055621233732447
$ws.Range("M23").Value = 1.066984593858212
$ws.Range("N23").Value = 1.051860250173579
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.044979157994918
$ws.Range("D24").Value = 1.052571911906387
$ws.Range("E24").Value = 1.053132854654594
$ws.Range("F24").Value = 1.06461137665127
$ws.Range("I24").Value = 1.048586019634171
$ws.Range("J24").Value = 1.051108926466662
$ws.Range("K24").Value = 1.055885143748495
$ws.Range("L24").Value = 1.056444197256919
$ws.Range("M24").Value = 1.067884539502571
$ws.Range("N24").Value = 1.052601621706123
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.046438530582853
$ws.Range("D25").Value = 1.053710091599884
$ws.Range("E25").Value = 1.054405386983533
$ws.Range("F25").Value = 1.065969737128236
$ws.Range("I25").Value = 1.049001826545232
$ws.Range("J25").Value = 1.051968035888574
$ws.Range("K25").Value = 1.056707378268772
$ws.Range("L25").Value = 1.05740057032189
$ws.Range("M25").Value = 1.068930375351783
$ws.Range("N25").Value = 1.053461951161956
